# repull data, push all data, mean calculation
# Update the dSF (F) column values for each game row to reflect the
# re-pulled data (delta between the starting pitcher's score and the
# score when they left the game, recalculated after pushing all data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = -1
    5  = -5
    6  = 1
    7  = 3
    8  = -1
    10 = -3
    11 = -5
    12 = -1
    13 = -3
    14 = -2
    15 = -2
    16 = -3
    18 = -1
    19 = -5
    20 = -1
    21 = -8
    22 = 4
    23 = -2
    24 = 4
    26 = -1
    27 = -8
    28 = 4
    29 = 1
    30 = 4
    31 = -3
    32 = 4
    33 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
